{"js": "// Replace the date line and all 25 two-digit-by-two-digit multiplication\n// answers in the practice sheet, in place, preserving run formatting.\n//\n// Each old value is unique within the document, so a plain text search +\n// InsertLocation.replace on the hit range swaps the text without touching\n// the surrounding <w:rPr> (font/size) already on that run.\nconst replacements = [\n  [\"2024-09-10 Tuesday\", \"2024-09-11 Wednesday\"],\n  [\"42\u00d799=4158\", \"94\u00d715=1410\"],\n  [\"36\u00d756=2016\", \"25\u00d735=875\"],\n  [\"97\u00d716=1552\", \"83\u00d736=2988\"],\n  [\"16\u00d793=1488\", \"89\u00d736=3204\"],\n  [\"89\u00d797=8633\", \"78\u00d767=5226\"],\n  [\"92\u00d711=1012\", \"23\u00d725=575\"],\n  [\"41\u00d712=492\", \"15\u00d720=300\"],\n  [\"38\u00d773=2774\", \"52\u00d737=1924\"],\n  [\"38\u00d741=1558\", \"56\u00d756=3136\"],\n  [\"28\u00d765=1820\", \"96\u00d780=7680\"],\n  [\"16\u00d772=1152\", \"47\u00d795=4465\"],\n  [\"12\u00d749=588\", \"28\u00d789=2492\"],\n  [\"93\u00d729=2697\", \"73\u00d733=2409\"],\n  [\"81\u00d795=7695\", \"25\u00d752=1300\"],\n  [\"13\u00d757=741\", \"27\u00d759=1593\"],\n  [\"33\u00d784=2772\", \"83\u00d781=6723\"],\n  [\"93\u00d727=2511\", \"59\u00d765=3835\"],\n  [\"56\u00d767=3752\", \"15\u00d771=1065\"],\n  [\"70\u00d724=1680\", \"75\u00d724=1800\"],\n  [\"92\u00d750=4600\", \"70\u00d722=1540\"],\n  [\"96\u00d784=8064\", \"19\u00d752=988\"],\n  [\"50\u00d738=1900\", \"67\u00d722=1474\"],\n  [\"60\u00d771=4260\", \"51\u00d724=1224\"],\n  [\"43\u00d739=1677\", \"29\u00d781=2349\"],\n  [\"18\u00d716=288\", \"80\u00d760=4800\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and all 25 two-digit-by-two-digit multiplication\n# answers in the practice sheet, in place, preserving run formatting.\n#\n# Each old value is unique in the document, so Find/Replace (wildcards off,\n# whole text match) on $d.Content swaps just the run text and leaves the\n# existing <w:rPr> (font/size) on that run untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-09-10 Tuesday\", \"2024-09-11 Wednesday\"),\n  @(\"42\u00d799=4158\", \"94\u00d715=1410\"),\n  @(\"36\u00d756=2016\", \"25\u00d735=875\"),\n  @(\"97\u00d716=1552\", \"83\u00d736=2988\"),\n  @(\"16\u00d793=1488\", \"89\u00d736=3204\"),\n  @(\"89\u00d797=8633\", \"78\u00d767=5226\"),\n  @(\"92\u00d711=1012\", \"23\u00d725=575\"),\n  @(\"41\u00d712=492\", \"15\u00d720=300\"),\n  @(\"38\u00d773=2774\", \"52\u00d737=1924\"),\n  @(\"38\u00d741=1558\", \"56\u00d756=3136\"),\n  @(\"28\u00d765=1820\", \"96\u00d780=7680\"),\n  @(\"16\u00d772=1152\", \"47\u00d795=4465\"),\n  @(\"12\u00d749=588\", \"28\u00d789=2492\"),\n  @(\"93\u00d729=2697\", \"73\u00d733=2409\"),\n  @(\"81\u00d795=7695\", \"25\u00d752=1300\"),\n  @(\"13\u00d757=741\", \"27\u00d759=1593\"),\n  @(\"33\u00d784=2772\", \"83\u00d781=6723\"),\n  @(\"93\u00d727=2511\", \"59\u00d765=3835\"),\n  @(\"56\u00d767=3752\", \"15\u00d771=1065\"),\n  @(\"70\u00d724=1680\", \"75\u00d724=1800\"),\n  @(\"92\u00d750=4600\", \"70\u00d722=1540\"),\n  @(\"96\u00d784=8064\", \"19\u00d752=988\"),\n  @(\"50\u00d738=1900\", \"67\u00d722=1474\"),\n  @(\"60\u00d771=4260\", \"51\u00d724=1224\"),\n  @(\"43\u00d739=1677\", \"29\u00d781=2349\"),\n  @(\"18\u00d716=288\", \"80\u00d760=4800\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n\n"}
